$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.4
$ws.Range("H2").Value = 2.75
$ws.Range("O2").Value = 1.73
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 2.48
$ws.Range("R2").Value = 1.55
$ws.Range("U2").Value = 6.2
$ws.Range("V2").Value = 1.13
$ws.Range("W2").Value = 8
$ws.Range("X2").Value = 1.08
$ws.Range("G3").Value = 2.2
$ws.Range("H3").Value = 2.85
$ws.Range("I3").Value = 3.8
$ws.Range("I4").Value = 2.35
$ws.Range("G9").Value = 1.62
$ws.Range("I9").Value = 5
$ws.Range("J9").Value = 2.2
$ws.Range("AA9").Value = 1.62
$ws.Range("AB9").Value = 2.2
$ws.Range("AK9").Value = 13
$ws.Range("AQ9").Value = 34
$ws.Range("AR9").Value = 34
$ws.Range("G16").Value = 5.1
$ws.Range("H16").Value = 4.1
$ws.Range("I16").Value = 1.57
$ws.Range("J16").Value = 5.1
$ws.Range("L16").Value = 2.1
$ws.Range("O16").Value = 1.22
$ws.Range("P16").Value = 3.85
$ws.Range("S16").Value = 1.55
$ws.Range("T16").Value = 2.1
$ws.Range("W16").Value = 2.6
$ws.Range("X16").Value = 1.44
$ws.Range("Z16").Value = 3.2
$ws.Range("AA16").Value = 1.65
$ws.Range("AB16").Value = 2.12
$ws.Range("AC16").Value = 15.5
$ws.Range("AE16").Value = 16
$ws.Range("AG16").Value = 45
$ws.Range("AH16").Value = 45
$ws.Range("AJ16").Value = 8
$ws.Range("AM16").Value = 7.8
$ws.Range("AN16").Value = 7.8
$ws.Range("AP16").Value = 11.5
$ws.Range("AQ16").Value = 11.75
$ws.Range("M17").Value = 1.11
$ws.Range("N17").Value = 6.5
$ws.Range("AB17").Value = 1.67
